$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.731.93"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "2.227.15"
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("E4").Value = "  +0.12%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "271.07"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +4.94%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "89.42"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +11.78%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.620"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -0.77%  "
$ws.Range("E8").Value = "  +0.03%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.604"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.36%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "45.91"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +6.02%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0917"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -0.96%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "7.76"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +9.67%  "
$ws.Range("E13").Value = "  +1.17%  "
$ws.Range("D14").Value = "2.562.87"
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("E15").Value = "  +2.03%  "
$ws.Range("D16").Value = "2.225.57"
$ws.Range("E16").Value = "  -0.36%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.790"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "43.665.06"
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "70.31"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("E21").Value = "  -1.26%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "2.34"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -0.21%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "231.99"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.79%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "8.57"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -8.85%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.08%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.51"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +12.66%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "10.91"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("E28").Value = "  +5.51%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.27"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +2.75%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "38.53"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -5.01%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "172.56"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +0.02%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.0907"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +1.95%  "
$ws.Range("E33").Value = "  +0.48%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "5.34"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +0.70%  "
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("E36").Value = "  -2.23%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.0351"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -4.45%  "
$ws.Range("E38").Value = "  -6.07%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "3.43"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +15.32%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "2.15"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +0.49%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "12.33"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -6.44%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.212"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +4.84%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "63.15"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +0.28%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "5.37"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -2.78%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "8.48"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -0.44%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.0986"
$cell.Style = "Normal"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "99.90"
$cell.Style = "Normal"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.15"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +1.74%  "
$ws.Range("E49").Value = "  +2.43%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.434"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -3.08%  "
$ws.Range("E51").Value = "  -5.26%  "
